$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.092.34"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.596.60"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.42"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.90"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.88"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.37"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0839"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.995.98"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.600.43"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.914"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.173.92"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  -4.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "285.19"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +10.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.36"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.90"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +6.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.78"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.20"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.65"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.80"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("E39").Value = "  +4.88%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.90"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +17.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.71"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.57"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.95"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.092.72"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.26"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.61"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +5.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.61"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -0.71%  "
